# Apply the edits described by the commit diff:
#  1. Rename sheet "LKT 8Beta3" -> "LKT 8HED3" (first sheet, rId1 / sheetId 4).
#  2. Fix a typo in the HED-tag shared string used by cell E4 of that sheet:
#     "Experimental-participant" -> "Experiment-participant".
#  3. Move the sheet's active-cell selection from E4 to E29.

$wb = $excel.ActiveWorkbook

# 1. Rename the first worksheet.
$ws = $wb.Worksheets.Item(1)
$ws.Name = "LKT 8HED3"

# 2. Correct the text typo inside the HED tag string stored in E4.
$ws.Range("E4").Value = "Agent-action, Participant-response, Correction, ((Human-agent, Experiment-participant), (Modify, (Car, Direction)))"

# 3. Update the active selection on the sheet from E4 to E29.
$ws.Activate()
$ws.Range("E29").Select()
